{"js": "const body = context.document.body;\n\n// --- Change 1 ---------------------------------------------------------\n// \"...execution of test scripts with a custom tool that utilizes Selenium\n// WebDriver.\" becomes \"...execution of 18% of regression scripts, along\n// with 37 change request scripts, with a custom tool that utilizes\n// Selenium WebDriver.\"\n\n// Replace \"test scripts\" with \"18% of regression scripts\" (keeps the\n// trailing \" scripts\" text identical, mirrors just swapping \"test\" out).\nconst testScriptsResults = body.search(\"test scripts\", { matchCase: true });\ntestScriptsResults.load(\"items\");\nawait context.sync();\nif (testScriptsResults.items.length === 0) {\n  throw new Error('Could not find \"test scripts\" to update.');\n}\ntestScriptsResults.items[0].insertText(\n  \"18% of regression scripts\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Insert the new \"along with 37 change request scripts\" clause right\n// before \"with a custom tool that utilizes Selenium WebDriver\".\nconst toolClauseResults = body.search(\n  \" with a custom tool that utilizes Selenium WebDriver\",\n  { matchCase: true }\n);\ntoolClauseResults.load(\"items\");\nawait context.sync();\nif (toolClauseResults.items.length === 0) {\n  throw new Error('Could not find the \"with a custom tool...\" clause to update.');\n}\ntoolClauseResults.items[0].insertText(\n  \", along with 37 change request scripts, with a custom tool that utilizes Selenium WebDriver\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Change 2 ---------------------------------------------------------\n// \"...completion of a test software, written in Python, that receives...\"\n// becomes \"...completion of test software, written in Python, that\n// receives...\" (drop the indefinite article \"a \").\nconst testSoftwareResults = body.search(\n  \"a test software, written in Python\",\n  { matchCase: true }\n);\ntestSoftwareResults.load(\"items\");\nawait context.sync();\nif (testSoftwareResults.items.length === 0) {\n  throw new Error('Could not find \"a test software, written in Python\" to update.');\n}\ntestSoftwareResults.items[0].insertText(\n  \"test software, written in Python\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$doc = $word.ActiveDocument\n\n# --- Change 1 -----------------------------------------------------------\n# \"...execution of test scripts with a custom tool that utilizes Selenium\n# WebDriver.\" becomes \"...execution of 18% of regression scripts, along\n# with 37 change request scripts, with a custom tool that utilizes\n# Selenium WebDriver.\"\n\n# 1a: \"test scripts\" -> \"18% of regression scripts\"\n$find1 = $doc.Content.Find\n$find1.ClearFormatting()\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.Text = \"test scripts\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"18% of regression scripts\"\n$find1.Execute($find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 1)\n\n# 1b: insert the new \"along with 37 change request scripts\" clause right\n# before \"with a custom tool that utilizes Selenium WebDriver\"\n$find2 = $doc.Content.Find\n$find2.ClearFormatting()\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.Text = \" with a custom tool that utilizes Selenium WebDriver\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \", along with 37 change request scripts, with a custom tool that utilizes Selenium WebDriver\"\n$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 1)\n\n# --- Change 2 -------------------------------------------------------------\n# \"...completion of a test software, written in Python, that receives...\"\n# becomes \"...completion of test software, written in Python, that\n# receives...\" (drop the indefinite article \"a \").\n$find3 = $doc.Content.Find\n$find3.ClearFormatting()\n$find3.MatchCase = $true\n$find3.MatchWholeWord = $false\n$find3.Text = \"a test software, written in Python\"\n$find3.Replacement.ClearFormatting()\n$find3.Replacement.Text = \"test software, written in Python\"\n$find3.Execute($find3.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 1)\n"}
